# Update "想去人数" (column F) figures across the sheets to match the
# newly generated output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 10739
$ws1.Range("F4").Value = 239
$ws1.Range("F6").Value = 7165
$ws1.Range("F7").Value = 144
$ws1.Range("F8").Value = 12819
$ws1.Range("F9").Value = 13249
$ws1.Range("F10").Value = 1333
$ws1.Range("F12").Value = 5581
$ws1.Range("F13").Value = 933
$ws1.Range("F20").Value = 1076
$ws1.Range("F21").Value = 1630
$ws1.Range("F26").Value = 3080
$ws1.Range("F28").Value = 2134
$ws1.Range("F31").Value = 1713
$ws1.Range("F32").Value = 1020
$ws1.Range("F33").Value = 867
$ws1.Range("F34").Value = 72
$ws1.Range("F36").Value = 3845
$ws1.Range("F37").Value = 4501
$ws1.Range("F42").Value = 3169
$ws1.Range("F47").Value = 46
$ws1.Range("F48").Value = 4324
$ws1.Range("F49").Value = 209

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 40
$ws2.Range("F23").Value = 80

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6795
$ws3.Range("F3").Value = 122
$ws3.Range("F4").Value = 280

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 10739
$ws4.Range("F4").Value = 7165
$ws4.Range("F5").Value = 144
$ws4.Range("F6").Value = 122
$ws4.Range("F7").Value = 12819
$ws4.Range("F8").Value = 13249
$ws4.Range("F10").Value = 1333
$ws4.Range("F12").Value = 5581
$ws4.Range("F13").Value = 933
$ws4.Range("F19").Value = 1076
$ws4.Range("F20").Value = 1630
$ws4.Range("F24").Value = 3080
$ws4.Range("F27").Value = 2134
$ws4.Range("F31").Value = 1713
$ws4.Range("F33").Value = 1020
$ws4.Range("F34").Value = 867
$ws4.Range("F35").Value = 72
$ws4.Range("F36").Value = 3845
$ws4.Range("F37").Value = 4501
$ws4.Range("F43").Value = 3169
$ws4.Range("F47").Value = 46
$ws4.Range("F48").Value = 4324
$ws4.Range("F49").Value = 209
